{"js": "// Update the worked-out binary-to-decimal conversion examples in the\n// second column of the table (cell 1,1 -> 0-indexed row 0, column 1).\n//\n// 1) \"1*1+1*2+1*4+1*8 = 15\"\n//      -> \"1111 = 1*8 + 1*4 + 1*2 + 1*1 = 15\"\n// 2) \"10 1010 = 1*1+0*2 + 1*1+0*2+1*4+0*8 = 6\"\n//      -> \"101010 =  \" + \"1*32 + 0*16 + 1*8 + 0*4 + 1*2 + 0*1 = 42\"\n// 3) \"111 1110 0101 = 1*1+1*2+1*4 + 1*1+1*2+1*4+0*8 + 0*1+1*2+0*4+1*8 = 24\"\n//      -> \"11111100101 = \" + \"1*1 + 0*2 + 1*4 + 0*8 + 0*16 + 1*32 + 1*64 + 1*128 + 1*256 + 1*512 + 1*1024 = 2021\"\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst cell = table.getCell(0, 1); // second column (\"Omrekenen\" answers)\n\nconst paragraphs = cell.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 1: single run -------------------------------------------\nconst p0 = paragraphs.items[0];\np0.insertText(\"1111 = 1*8 + 1*4 + 1*2 + 1*1 = 15\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Paragraph 2: two runs (\"101010 =  \" / \"1*32 + ... = 42\") ----------\nconst p1 = paragraphs.items[1];\np1.insertText(\"101010 =  \", Word.InsertLocation.replace);\nawait context.sync();\n\nconst p1Run2 = p1.insertText(\n  \"1*32 + 0*16 + 1*8 + 0*4 + 1*2 + 0*1 = 42\",\n  Word.InsertLocation.end\n);\n// Touch (and revert) a formatting property so the new text is emitted as\n// its own run rather than being merged back into the previous one.\np1Run2.font.bold = true;\nawait context.sync();\np1Run2.font.bold = false;\nawait context.sync();\n\n// --- Paragraph 3: two runs (\"11111100101 = \" / \"1*1 + ... = 2021\") -----\nconst p2 = paragraphs.items[2];\np2.insertText(\"11111100101 = \", Word.InsertLocation.replace);\nawait context.sync();\n\nconst p2Run2 = p2.insertText(\n  \"1*1 + 0*2 + 1*4 + 0*8 + 0*16 + 1*32 + 1*64 + 1*128 + 1*256 + 1*512 + 1*1024 = 2021\",\n  Word.InsertLocation.end\n);\np2Run2.font.bold = true;\nawait context.sync();\np2Run2.font.bold = false;\nawait context.sync();\n", "ps1": "# Update the worked-out binary-to-decimal conversion examples in the\n# second column of the table (\"Omrekenen\" row), using Find & Replace so\n# the edit is anchored on content rather than brittle indices.\n#\n# 1) \"1*1+1*2+1*4+1*8 = 15\"\n#      -> \"1111 = 1*8 + 1*4 + 1*2 + 1*1 = 15\"\n# 2) \"10 1010 = 1*1+0*2 + 1*1+0*2+1*4+0*8 = 6\"\n#      -> \"101010 =  \" + \"1*32 + 0*16 + 1*8 + 0*4 + 1*2 + 0*1 = 42\"\n# 3) \"111 1110 0101 = 1*1+1*2+1*4 + 1*1+1*2+1*4+0*8 + 0*1+1*2+0*4+1*8 = 24\"\n#      -> \"11111100101 = \" + \"1*1 + 0*2 + 1*4 + 0*8 + 0*16 + 1*32 + 1*64 + 1*128 + 1*256 + 1*512 + 1*1024 = 2021\"\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: single run, straight replace ---------------------------\n$rng = $d.Content\n$null = $rng.Find.Execute(\n    \"1*1+1*2+1*4+1*8 = 15\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"1111 = 1*8 + 1*4 + 1*2 + 1*1 = 15\", 2)\n\n# --- Paragraph 2: replace with first run's text, then append a 2nd run ---\n$d = $word.ActiveDocument\n$rng = $d.Content\n$null = $rng.Find.Execute(\n    \"10 1010 = 1*1+0*2 + 1*1+0*2+1*4+0*8 = 6\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"101010 =  \", 2)\n\n$d = $word.ActiveDocument\n$rng2 = $d.Content\n$null = $rng2.Find.Execute(\"101010 =  \")\n$rng2.Collapse(0)\n$rng2.InsertAfter(\"1*32 + 0*16 + 1*8 + 0*4 + 1*2 + 0*1 = 42\")\n# Touch (and revert) a formatting property so the appended text lands in\n# its own run instead of merging back into the preceding one.\n$rng2.Font.Bold = 1\n$rng2.Font.Bold = 0\n\n# --- Paragraph 3: replace with first run's text, then append a 2nd run ---\n$d = $word.ActiveDocument\n$rng = $d.Content\n$null = $rng.Find.Execute(\n    \"111 1110 0101 = 1*1+1*2+1*4 + 1*1+1*2+1*4+0*8 + 0*1+1*2+0*4+1*8 = 24\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"11111100101 = \", 2)\n\n$d = $word.ActiveDocument\n$rng3 = $d.Content\n$null = $rng3.Find.Execute(\"11111100101 = \")\n$rng3.Collapse(0)\n$rng3.InsertAfter(\"1*1 + 0*2 + 1*4 + 0*8 + 0*16 + 1*32 + 1*64 + 1*128 + 1*256 + 1*512 + 1*1024 = 2021\")\n$rng3.Font.Bold = 1\n$rng3.Font.Bold = 0\n"}
